$wb = $excel.ActiveWorkbook

# "WithdrawalAccount" is the sheet being updated with a new header row
# (UserName / Password) above the existing credentials row.
$ws = $wb.Worksheets.Item("WithdrawalAccount")

# Insert a new row above the existing data row so the current A1:B1
# content shifts down to A2:B2.
$ws.Rows.Item(1).Insert() | Out-Null

$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Match the style (left/top aligned) used for the analogous header row
# on the "User" sheet.
$headerSrc = $wb.Worksheets.Item("User").Range("A1:B1")
$headerSrc.Copy()
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$userSheet = $wb.Worksheets.Item("User")
$userSheet.Range("B10").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("B6").Select() | Out-Null
